# Applies the RDIT document edits described by the target diff:
#  1. Resize the first three table rows' cell widths (4788 -> 4681 / 4669 twips)
#     and delete the trailing empty row.
#  2. Clean up the two "The user is able to ..." sentences that had
#     proofErr-wrapped "is able to" runs, merging them into single runs.
#  3. Insert two new sub-bullets in the Design/Implementation sections.
#  4. Insert three new testing steps.
#  5. Drop the stray "_GoBack" bookmark pair left in the final paragraph.
#
# NOTE: once a Table/Row/Cell object has been touched, `$d.Paragraphs.Item(n)`
# becomes unreliable in this host (it keeps returning the same paragraph).
# `$d.Content.Paragraphs` does not have that problem, so it is used
# throughout instead.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Table cell widths + row removal
# ---------------------------------------------------------------------
$table = $d.Tables.Item(1)

for ($r = 1; $r -le 3; $r++) {
    $row = $table.Rows.Item($r)
    $row.Cells.Item(1).Width = 234.05   # 4681 twips
    $row.Cells.Item(2).Width = 233.45   # 4669 twips
}

# Remove the trailing blank row (originally row 4).
if ($table.Rows.Count -ge 4) {
    $table.Rows.Item(4).Delete()
}

# ---------------------------------------------------------------------
# 2. Merge the proofErr-split "The user is able to ..." runs
# ---------------------------------------------------------------------
$d.Content.Find.Execute("The user is able to type", $true, $false, $false, $false, $false, $true, 1, $false, "The user is able to type", 2) | Out-Null
$d.Content.Find.Execute("The user is able to browse for a", $true, $false, $false, $false, $false, $true, 1, $false, "The user is able to browse for a", 2) | Out-Null

# ---------------------------------------------------------------------
# Helper: find a paragraph by its exact (trimmed) text.
# ---------------------------------------------------------------------
function Get-ParagraphByText($text) {
    $paras = $d.Content.Paragraphs
    $count = $paras.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 3. New sub-bullets under Design / Implementation
# ---------------------------------------------------------------------
$anchor = Get-ParagraphByText "Add a QTextEdit that has an example file path as placeholder text"
$anchor.Range.InsertParagraphAfter()
$newPara = $anchor.Next()
$newPara.Range.Text = "If an invalid file path is specified then upon processing a QMessageBox will prompt the user with an error message"
$newPara.Range.ListFormat.ListLevelNumber = 3

$anchor2 = Get-ParagraphByText "Upon clicking the QButton prompt the user with an open QFileDialog"
$anchor2.Range.InsertParagraphAfter()
$newPara2 = $anchor2.Next()
$newPara2.Range.Text = "If the QFileDialog is cancelled, nothing populates the input file QTextEdit entry field"
$newPara2.Range.ListFormat.ListLevelNumber = 3

# ---------------------------------------------------------------------
# 4. New testing steps
# ---------------------------------------------------------------------
$quote1 = [char]0x201c
$quote2 = [char]0x201d

$testAnchor = Get-ParagraphByText "Verify that a button with a folder icon exists to the right of the text entry field"
$testAnchor.Range.InsertParagraphAfter()
$t1 = $testAnchor.Next()
$t1.Range.Text = "Enter an invalid path into the input file text entry field. For Example: " + $quote1 + "3>N>N>N" + $quote2
$t1.Range.ListFormat.ListLevelNumber = 2

$t1.Range.InsertParagraphAfter()
$t2 = $t1.Next()
$t2.Range.Text = "Click the " + $quote1 + "Process" + $quote2 + " button and verify that a message box prompts the user with an error message, the border of the text entry field turns red, and no processing happens"
$t2.Range.ListFormat.ListLevelNumber = 2

$t2.Range.InsertParagraphAfter()
$t3 = $t2.Next()
$t3.Range.Text = "Repeat 4.5 and 4.6 but clear the entry field"
$t3.Range.ListFormat.ListLevelNumber = 2

# ---------------------------------------------------------------------
# 5. Bookmarks: drop the stray "_GoBack" pair inserted in the final
#    paragraph (the remaining "_Hlk46782967" bookmark is renumbered
#    automatically by Word once "_GoBack" is gone).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
